$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 44665
$ws.Range("J3").Value = 44665
$ws.Range("L3").Value = 44665
$ws.Range("N3").Value = -44893

$ws.Range("H88").Value = 4761.36
$ws.Range("J88").Value = 4212.5264
$ws.Range("L88").Value = 4212.5264
$ws.Range("N88").Value = -5024.5264

$ws.Range("H91").Value = 4761.36
$ws.Range("J91").Value = 4212.5264
$ws.Range("L91").Value = 4212.5264
$ws.Range("N91").Value = -7020.5264

$ws.Range("H96").Value = 1371.9286
$ws.Range("I96").Value = 1578.5
$ws.Range("J96").Value = 1096.5
$ws.Range("K96").Value = 4735.5
$ws.Range("L96").Value = 3289.5
$ws.Range("M96").Value = -3362.5
$ws.Range("N96").Value = -6035.5

$ws.Range("H102").Value = 44665
$ws.Range("J102").Value = 44665
$ws.Range("L102").Value = 44665
$ws.Range("N102").Value = -51155

$ws.Range("H107").Value = 1478.1666
$ws.Range("I107").Value = 1476.1818
$ws.Range("K107").Value = 1476.1818
$ws.Range("M107").Value = 443.8181999999999

$ws.Range("H111").Value = 2181.2307
$ws.Range("I111").Value = 1882.4
$ws.Range("J111").Value = 3177.3333
$ws.Range("K111").Value = 5647.200000000001
$ws.Range("L111").Value = 9531.999899999999
$ws.Range("M111").Value = -2580.200000000001
$ws.Range("N111").Value = -15665.9999

$ws.Range("H135").Value = 9625.788
$ws.Range("I135").Value = 6082.55
$ws.Range("J135").Value = 15076.923
$ws.Range("K135").Value = 54742.95
$ws.Range("L135").Value = 135692.307
$ws.Range("M135").Value = -52207.95
$ws.Range("N135").Value = -140762.307

$ws.Range("H138").Value = 4167.9
$ws.Range("J138").Value = 4930.9165
$ws.Range("L138").Value = 14792.7495
$ws.Range("N138").Value = -25072.7495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16395467
$ws.Range("I32").Value = 17545926
$ws.Range("J32").Value = 1424.25
$ws.Range("K32").Value = 17545926
$ws.Range("L32").Value = 1424.25
$ws.Range("M32").Value = -17545639
$ws.Range("N32").Value = -1998.25

$ws.Range("H132").Value = 1731.8948
$ws.Range("I132").Value = 1570.8572
$ws.Range("J132").Value = 2182.8
$ws.Range("K132").Value = 4712.571599999999
$ws.Range("L132").Value = 6548.400000000001
$ws.Range("M132").Value = -2182.571599999999
$ws.Range("N132").Value = -11608.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4768.6665
$ws.Range("J105").Value = 10011
$ws.Range("L105").Value = 10011
$ws.Range("N105").Value = -13505

$ws.Range("H107").Value = 8001
$ws.Range("I107").Value = 8418.5
$ws.Range("K107").Value = 8418.5
$ws.Range("M107").Value = -6498.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5586.75
$ws.Range("I16").Value = 8175.4287
$ws.Range("J16").Value = 1962.6
$ws.Range("K16").Value = 8175.4287
$ws.Range("L16").Value = 1962.6
$ws.Range("M16").Value = -7888.4287
$ws.Range("N16").Value = -2536.6

$ws.Range("H86").Value = 26329
$ws.Range("I86").Value = 29406.715
$ws.Range("J86").Value = 20173.572
$ws.Range("K86").Value = 29406.715
$ws.Range("L86").Value = 20173.572
$ws.Range("M86").Value = -28283.715
$ws.Range("N86").Value = -22419.572

$ws.Range("H89").Value = 26329
$ws.Range("I89").Value = 29406.715
$ws.Range("J89").Value = 20173.572
$ws.Range("K89").Value = 147033.575
$ws.Range("L89").Value = 100867.86
$ws.Range("M89").Value = -141417.575
$ws.Range("N89").Value = -112099.86

$ws.Range("H113").Value = 5586.75
$ws.Range("I113").Value = 8175.4287
$ws.Range("J113").Value = 1962.6
$ws.Range("K113").Value = 8175.4287
$ws.Range("L113").Value = 1962.6
$ws.Range("M113").Value = -6005.4287
$ws.Range("N113").Value = -6302.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 512.0833
$ws.Range("I5").Value = 308.75
$ws.Range("J5").Value = 613.75
$ws.Range("K5").Value = 926.25
$ws.Range("L5").Value = 1841.25
$ws.Range("M5").Value = -814.25
$ws.Range("N5").Value = -2065.25

$ws.Range("H11").Value = 120
$ws.Range("I11").Value = 120
$ws.Range("K11").Value = 360
$ws.Range("M11").Value = -220

$ws.Range("H55").Value = 437735.12
$ws.Range("J55").Value = 503180.2
$ws.Range("L55").Value = 1509540.6
$ws.Range("N55").Value = -1509894.6

$ws.Range("H68").Value = 1785.1428
$ws.Range("I68").Value = 998.5
$ws.Range("J68").Value = 2099.8
$ws.Range("K68").Value = 2995.5
$ws.Range("L68").Value = 6299.400000000001
$ws.Range("M68").Value = -2184.5
$ws.Range("N68").Value = -7921.400000000001

$ws.Range("H71").Value = 1785.1428
$ws.Range("I71").Value = 998.5
$ws.Range("J71").Value = 2099.8
$ws.Range("K71").Value = 8986.5
$ws.Range("L71").Value = 18898.2
$ws.Range("M71").Value = -4930.5
$ws.Range("N71").Value = -27010.2

$ws.Range("H88").Value = 4821.357
$ws.Range("I88").Value = 3900
$ws.Range("J88").Value = 5333.222
$ws.Range("K88").Value = 11700
$ws.Range("L88").Value = 15999.666
$ws.Range("M88").Value = -11272
$ws.Range("N88").Value = -16855.666

$ws.Range("H91").Value = 4821.357
$ws.Range("I91").Value = 3900
$ws.Range("J91").Value = 5333.222
$ws.Range("K91").Value = 11700
$ws.Range("L91").Value = 15999.666
$ws.Range("M91").Value = -10218
$ws.Range("N91").Value = -18963.666

$ws.Range("H107").Value = 1259.5385
$ws.Range("J107").Value = 1240
$ws.Range("L107").Value = 3720
$ws.Range("N107").Value = -7560

$ws.Range("H113").Value = 431.625
$ws.Range("J113").Value = 452.5
$ws.Range("L113").Value = 1357.5
$ws.Range("N113").Value = -5697.5

$ws.Range("H122").Value = 5476.5557
$ws.Range("I122").Value = 738.2857
$ws.Range("J122").Value = 8491.817999999999
$ws.Range("K122").Value = 6644.571300000001
$ws.Range("L122").Value = 76426.36199999999
$ws.Range("M122").Value = -4194.571300000001
$ws.Range("N122").Value = -81326.36199999999

$ws.Range("H135").Value = 512.0833
$ws.Range("I135").Value = 308.75
$ws.Range("J135").Value = 613.75
$ws.Range("K135").Value = 2778.75
$ws.Range("L135").Value = 5523.75
$ws.Range("M135").Value = -243.75
$ws.Range("N135").Value = -10593.75

$ws.Range("H140").Value = 1044.7142
$ws.Range("I140").Value = 932.7692
$ws.Range("J140").Value = 2500
$ws.Range("K140").Value = 2798.3076
$ws.Range("L140").Value = 7500
$ws.Range("M140").Value = 2381.6924
$ws.Range("N140").Value = -17860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2825473.5
$ws.Range("I24").Value = 7500666.5
$ws.Range("J24").Value = 20357.6
$ws.Range("K24").Value = 7500666.5
$ws.Range("L24").Value = 20357.6
$ws.Range("M24").Value = -7500493.5
$ws.Range("N24").Value = -20703.6

$ws.Range("H132").Value = 7288.579
$ws.Range("I132").Value = 6838.4
$ws.Range("K132").Value = 20515.2
$ws.Range("M132").Value = -17985.2

$ws.Range("H135").Value = 63556
$ws.Range("J135").Value = 63556
$ws.Range("L135").Value = 63556
$ws.Range("N135").Value = -73696

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H132").Value = 7199.2646
$ws.Range("I132").Value = 7099.0347
$ws.Range("J132").Value = 7780.6
$ws.Range("K132").Value = 21297.1041
$ws.Range("L132").Value = 23341.8
$ws.Range("M132").Value = -18767.1041
$ws.Range("N132").Value = -28401.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -760

$ws.Range("H29").Value = 8336.666999999999
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 8336.666999999999
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 8336.666999999999
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -8916.666999999999

$ws.Range("H81").Value = 7966530
$ws.Range("I81").Value = 22031
$ws.Range("J81").Value = 27827778
$ws.Range("K81").Value = 44062
$ws.Range("L81").Value = 55655556
$ws.Range("M81").Value = -43001
$ws.Range("N81").Value = -55657678

$ws.Range("H84").Value = 7966530
$ws.Range("I84").Value = 22031
$ws.Range("J84").Value = 27827778
$ws.Range("K84").Value = 220310
$ws.Range("L84").Value = 278277780
$ws.Range("M84").Value = -215006
$ws.Range("N84").Value = -278288388

$ws.Range("H122").Value = 1651.4166
$ws.Range("I122").Value = 1352.5
$ws.Range("J122").Value = 2249.25
$ws.Range("K122").Value = 4057.5
$ws.Range("L122").Value = 6747.75
$ws.Range("M122").Value = -1607.5
$ws.Range("N122").Value = -11647.75

$ws.Range("H124").Value = 199833.33
$ws.Range("J124").Value = 199833.33
$ws.Range("L124").Value = 199833.33
$ws.Range("N124").Value = -209653.33

$ws.Range("H136").Value = 2830.1282
$ws.Range("I136").Value = 2497.6
$ws.Range("K136").Value = 7492.799999999999
$ws.Range("M136").Value = -4942.799999999999
